# Word COM-interop script: rewrite the "black holes" essay into a
# "chemistry" essay (title, author, email, and body paragraphs),
# and normalize the font name from the misspelled "TimesNewToman"
# to "Times New Roman".

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / byline -------------------------------------------------
Replace-Text "The Allure of Black Holes: Unveiling Cosmic Mysteries" "The Intricate Dance of Chemistry: Transforming the Invisible into the Tangible"
Replace-Text "Marcus Hawthorne" "Jennifer Carter"
Replace-Text "marcus" "jennifer"
Replace-Text "hawthorne@luminary" "carter@growingscholaredu"
Replace-Text "academia" "org"

# --- First body paragraph -------------------------------------------
Replace-Text "In the vast cosmic tapestry, black holes captivate like no other celestial phenomenon" "Chemistry: the study of matter and its transformations"
Replace-Text " Their enigmatic allure stems from a paradoxical fusion of immense gravitational pull and a profound absence, an interplay between the finite and the infinite" " Not mere equations on a paper, it is an intricate dance of atoms and molecules, a symphony of elements interacting"
Replace-Text " Black holes shroud themselves within an event horizon, a point of no return where time and space become distorted, swallowed by the inescapable force of gravity" " In this realm, the invisible becomes tangible, the abstract takes physical form"
Replace-Text " Yet, they also hold the key to unlocking some of the universe's most profound secrets, tantalizing scientists to seek a deeper understanding of these enigmatic cosmic entities" " We explore the building blocks of the universe, from the smallest subatomic particles to the vast array of compounds that make up our world"

Replace-Text "The study of black holes has evolved from a theoretical exploration to a dynamic, observational field" "Beneath the placid surface of everyday objects lies a swirling vortex of chemical reactions"
Replace-Text " From the groundbreaking work of Karl Schwarzschild and Albert Einstein to the recent Event Horizon Telescope project, scientific advancements have provided a wealth of empirical insights about these formidable cosmic marvels" " The rust on metal, the ripening of fruit, the flame of a candle: each phenomenon a chemical metamorphosis"
Replace-Text " They have been discovered across diverse cosmic scales, from stellar black holes formed through the collapse of massive stars to gargantuan supermassive black holes at the heart of galaxies, each influencing their surroundings in profound ways" " The properties of substances, their colors, textures, and reactivities, stem from the intricate arrangements of constituent elements"

# Insert two new sentences after "...constituent elements." (before the
# double line-break that precedes "Unveiling the secrets...")
$find = $d.Content.Find
$find.Execute("constituent elements.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $find.Parent
$ins.Collapse(0)
$ins.InsertAfter(" Discoveries in chemistry have profoundly shaped our lives, from the medicines we take to the foods we eat, and even the materials that construct our homes")
$ins.Font.Name = "TimesNewToman"
$ins.Font.Size = 12
$ins.Font.Color = 0
$ins.Collapse(0)
$ins.InsertAfter(".")
$ins.Font.Name = "TimesNewToman"
$ins.Font.Size = 12
$ins.Font.Color = 0

Replace-Text "Black holes act as cosmic laboratories, revealing extreme conditions and challenging fundamental physical theories" "Unveiling the secrets of chemistry grants us the power to manipulate matter, to reshape the world around us"
Replace-Text " They are a frontier in physics, prompting scientists to explore the nature of spacetime curvature, quantum gravity, and the eventual fate of our universe" " We synthesize new substances with desirable properties, devising innovative solutions to global challenges"
Replace-Text " Moreover, they offer a unique window into the evolution of stars, the formation and merger of galaxies, and the birth of gravitational waves, ripples in spacetime that carry valuable information about the cosmos' past" " The field stands at the forefront of modern scientific inquiry, offering answers to some of the universe's most fundamental questions"

# Insert a new sentence after "...fundamental questions." (before the
# final, unchanged, trailing "." run of this paragraph)
$find2 = $d.Content.Find
$find2.Execute("fundamental questions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins2 = $find2.Parent
$ins2.Collapse(0)
$ins2.InsertAfter(" It is a discipline that bridges the gap between theory and practice, seamlessly intertwining abstract concepts with tangible applications")
$ins2.Font.Name = "TimesNewToman"
$ins2.Font.Size = 12
$ins2.Font.Color = 0

# --- Summary paragraph ------------------------------------------------
Replace-Text "Black holes stand as captivating cosmic enigmas, revealing extreme conditions that challenge our comprehension of the universe" "Chemistry, the study of matter and its transformations, is an intricate dance of atoms and molecules"
Replace-Text " They are the ultimate tests of our physical theories, pushing the boundaries of scientific knowledge and prompting profound questions about the nature of spacetime, quantum gravity, and the fate of our cosmos" " It unveils the invisible, giving tangible form to the abstract"
Replace-Text " The study of black holes promises to unveil cosmic mysteries, shedding light on the birth of stars, the evolution of galaxies, and the gravitational waves that carry hidden tales of our universe's history" " Our understanding of chemistry has revolutionized our lives, impacting medicine, agriculture, and material science"
Replace-Text " These cosmic wonders continue to fuel our scientific curiosity, driving advancements in astrophysics and providing invaluable insights into the fabric of our universe" " As we delve deeper into this realm, we gain unprecedented control over matter, harnessing its power to address global challenges and expand our knowledge of the universe"

# --- Trailing empty paragraph ------------------------------------------
$d.Content.InsertParagraphAfter()

# --- Global font-name fix: TimesNewToman -> Times New Roman ------------
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

Write-Host "edit complete"
